$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 29000
$ws.Range("I21").Value = 30000
$ws.Range("J21").Value = 28750
$ws.Range("K21").Value = 30000
$ws.Range("L21").Value = 28750
$ws.Range("M21").Value = -29532
$ws.Range("N21").Value = -29686
$ws.Range("H23").Value = 29000
$ws.Range("I23").Value = 30000
$ws.Range("J23").Value = 28750
$ws.Range("K23").Value = 30000
$ws.Range("L23").Value = 28750
$ws.Range("M23").Value = -29766
$ws.Range("N23").Value = -29218
$ws.Range("H98").Value = 2749.625
$ws.Range("I98").Value = 1965.079
$ws.Range("J98").Value = 5730.9
$ws.Range("K98").Value = 1965.079
$ws.Range("L98").Value = 5730.9
$ws.Range("M98").Value = -467.079
$ws.Range("N98").Value = -8726.9
$ws.Range("H122").Value = 2749.625
$ws.Range("I122").Value = 1965.079
$ws.Range("J122").Value = 5730.9
$ws.Range("K122").Value = 5895.237
$ws.Range("L122").Value = 17192.7
$ws.Range("M122").Value = -3445.237
$ws.Range("N122").Value = -22092.7
$ws.Range("H135").Value = 859.0909
$ws.Range("I135").Value = 707.8946999999999
$ws.Range("K135").Value = 6371.052299999999
$ws.Range("M135").Value = -3836.052299999999
$ws.Range("H137").Value = 1237.2084
$ws.Range("I137").Value = 1055.1333
$ws.Range("J137").Value = 1540.6666
$ws.Range("K137").Value = 3165.3999
$ws.Range("L137").Value = 4621.9998
$ws.Range("M137").Value = -615.3998999999999
$ws.Range("N137").Value = -9721.9998
$ws.Range("H138").Value = 4304.7964
$ws.Range("I138").Value = 2246.3096
$ws.Range("J138").Value = 9390.471
$ws.Range("K138").Value = 6738.9288
$ws.Range("L138").Value = 28171.413
$ws.Range("M138").Value = -1598.9288
$ws.Range("N138").Value = -38451.413

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8553.362999999999
$ws.Range("I32").Value = 9457.469999999999
$ws.Range("J32").Value = 3128.7273
$ws.Range("K32").Value = 9457.469999999999
$ws.Range("L32").Value = 3128.7273
$ws.Range("M32").Value = -9170.469999999999
$ws.Range("N32").Value = -3702.7273
$ws.Range("H61").Value = 1563.1562
$ws.Range("I61").Value = 1439.65
$ws.Range("J61").Value = 1769
$ws.Range("K61").Value = 1439.65
$ws.Range("L61").Value = 1769
$ws.Range("M61").Value = -1227.65
$ws.Range("N61").Value = -2193
$ws.Range("H74").Value = 973.7879
$ws.Range("I74").Value = 956.26086
$ws.Range("J74").Value = 1014.1
$ws.Range("K74").Value = 956.26086
$ws.Range("L74").Value = 1014.1
$ws.Range("M74").Value = -82.26085999999998
$ws.Range("N74").Value = -2762.1
$ws.Range("H77").Value = 973.7879
$ws.Range("I77").Value = 956.26086
$ws.Range("J77").Value = 1014.1
$ws.Range("K77").Value = 4781.3043
$ws.Range("L77").Value = 5070.5
$ws.Range("M77").Value = -413.3042999999998
$ws.Range("N77").Value = -13806.5
$ws.Range("H132").Value = 1678.0167
$ws.Range("I132").Value = 1035.2903
$ws.Range("J132").Value = 2365.0688
$ws.Range("K132").Value = 3105.8709
$ws.Range("L132").Value = 7095.2064
$ws.Range("M132").Value = -575.8708999999999
$ws.Range("N132").Value = -12155.2064
$ws.Range("H136").Value = 1563.1562
$ws.Range("I136").Value = 1439.65
$ws.Range("J136").Value = 1769
$ws.Range("K136").Value = 4318.950000000001
$ws.Range("L136").Value = 5307
$ws.Range("M136").Value = -1768.950000000001
$ws.Range("N136").Value = -10407

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1875.7872
$ws.Range("I134").Value = 1607.2703
$ws.Range("J134").Value = 2869.3
$ws.Range("K134").Value = 4821.810899999999
$ws.Range("L134").Value = 8607.900000000001
$ws.Range("M134").Value = -2286.810899999999
$ws.Range("N134").Value = -13677.9

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1840.3572
$ws.Range("I31").Value = 1244.742
$ws.Range("J31").Value = 3518.9092
$ws.Range("K31").Value = 1244.742
$ws.Range("L31").Value = 3518.9092
$ws.Range("M31").Value = -949.742
$ws.Range("N31").Value = -4108.9092
$ws.Range("H34").Value = 1840.3572
$ws.Range("I34").Value = 1244.742
$ws.Range("J34").Value = 3518.9092
$ws.Range("K34").Value = 1244.742
$ws.Range("L34").Value = 3518.9092
$ws.Range("M34").Value = -1042.742
$ws.Range("N34").Value = -3922.9092
$ws.Range("H58").Value = 598474.0600000001
$ws.Range("I58").Value = 862237.25
$ws.Range("J58").Value = 1536.3158
$ws.Range("K58").Value = 862237.25
$ws.Range("L58").Value = 1536.3158
$ws.Range("M58").Value = -862034.25
$ws.Range("N58").Value = -1942.3158
$ws.Range("H132").Value = 266200.03
$ws.Range("I132").Value = 330556.16
$ws.Range("J132").Value = 2339.8
$ws.Range("K132").Value = 991668.48
$ws.Range("L132").Value = 7019.400000000001
$ws.Range("M132").Value = -989138.48
$ws.Range("N132").Value = -12079.4
$ws.Range("H134").Value = 1066.3385
$ws.Range("I134").Value = 852.25
$ws.Range("J134").Value = 1670.8235
$ws.Range("K134").Value = 2556.75
$ws.Range("L134").Value = 5012.470499999999
$ws.Range("M134").Value = -21.75
$ws.Range("N134").Value = -10082.4705
$ws.Range("H136").Value = 598474.0600000001
$ws.Range("I136").Value = 862237.25
$ws.Range("J136").Value = 1536.3158
$ws.Range("K136").Value = 2586711.75
$ws.Range("L136").Value = 4608.9474
$ws.Range("M136").Value = -2584161.75
$ws.Range("N136").Value = -9708.947400000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 27163.666
$ws.Range("J92").Value = 27163.666
$ws.Range("L92").Value = 27163.666
$ws.Range("N92").Value = -30907.666
$ws.Range("H132").Value = 1478.7778
$ws.Range("I132").Value = 912.3889
$ws.Range("K132").Value = 2737.1667
$ws.Range("M132").Value = -207.1667000000002

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2096
$ws.Range("I82").Value = 1291.4286
$ws.Range("J82").Value = 2900.5715
$ws.Range("K82").Value = 1291.4286
$ws.Range("L82").Value = 2900.5715
$ws.Range("M82").Value = -930.4286
$ws.Range("N82").Value = -3622.5715
$ws.Range("H85").Value = 2096
$ws.Range("I85").Value = 1291.4286
$ws.Range("J85").Value = 2900.5715
$ws.Range("K85").Value = 1291.4286
$ws.Range("L85").Value = 2900.5715
$ws.Range("M85").Value = -43.42859999999996
$ws.Range("N85").Value = -5396.5715
$ws.Range("H132").Value = 2332.3845
$ws.Range("I132").Value = 1848.2258
$ws.Range("J132").Value = 4208.5
$ws.Range("K132").Value = 5544.6774
$ws.Range("L132").Value = 12625.5
$ws.Range("M132").Value = -3014.6774
$ws.Range("N132").Value = -17685.5
$ws.Range("H136").Value = 2921.1096
$ws.Range("I136").Value = 3047.93
$ws.Range("J136").Value = 2469.3125
$ws.Range("K136").Value = 9143.789999999999
$ws.Range("L136").Value = 7407.9375
$ws.Range("M136").Value = -6593.789999999999
$ws.Range("N136").Value = -12507.9375

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1094.4565
$ws.Range("I132").Value = 766.6389
$ws.Range("J132").Value = 2274.6
$ws.Range("K132").Value = 2299.9167
$ws.Range("L132").Value = 6823.799999999999
$ws.Range("M132").Value = 230.0832999999998
$ws.Range("N132").Value = -11883.8
$ws.Range("H136").Value = 4478.5713
$ws.Range("I136").Value = 6416.6665
$ws.Range("J136").Value = 3025
$ws.Range("K136").Value = 19249.9995
$ws.Range("L136").Value = 9075
$ws.Range("M136").Value = -16699.9995
$ws.Range("N136").Value = -14175
